$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'white athletic leggings mens'
$ws.Range("A2").Value = 'supreme basketball tights'
$ws.Range("A3").Value = 'snowmobiling knee pads'
$ws.Range("A4").Value = 'photographer knee pads'
$ws.Range("A5").Value = 'raymens leggings'
$ws.Range("A6").Value = 'training basketball youth'
$ws.Range("A7").Value = 'baleaf men compression pants'
$ws.Range("A8").Value = 'winter leggings men'
$ws.Range("A9").Value = 'kids pants with knee pads'
$ws.Range("A10").Value = 'eclipse knee pads'
$ws.Range("A11").Value = 'pilates knee pads'
$ws.Range("A12").Value = 'fitted mens tights'
$ws.Range("A13").Value = 'mens knee pads bmx'
$ws.Range("A14").Value = 'canoe knee pad'
$ws.Range("A15").Value = 'graduated compression pants'
$ws.Range("A16").Value = 'graduated compression leggings'
$ws.Range("A17").Value = 'soft knee pads'
$ws.Range("A18").Value = 'knee pads nike'
$ws.Range("A19").Value = 'knee pad toddler'
$ws.Range("A20").Value = 'knee pad yoga'
$ws.Range("A21").Value = 'knee pads dance'
$ws.Range("A22").Value = 'knee pads mizuno'
$ws.Range("A23").Value = 'elbow knee pads'
$ws.Range("A24").Value = 'nee pads basketball'
$ws.Range("A25").Value = 'compression pants set men'
$ws.Range("A26").Value = 'mens basketball snap pants'
$ws.Range("A27").Value = 'nike compression pants for youth'
$ws.Range("A28").Value = 'mens compression tights cold weather'
$ws.Range("A29").Value = 'mens under armour basketball tights'
$ws.Range("A30").Value = 'mens tights pockets'
$ws.Range("A31").Value = 'rollerblades knee pads'
$ws.Range("A32").Value = 'knee pads 3xl'
$ws.Range("A33").Value = 'knee pads 8'
$ws.Range("A34").Value = 'knee pad hard'
$ws.Range("A35").Value = 'men leggings fleece'
$ws.Range("A36").Value = 'teflex knee pads'
$ws.Range("A37").Value = 'man leggings thermal'
$ws.Range("A38").Value = 'sailing knee pad'
$ws.Range("A39").Value = 'knee pads downhill'
$ws.Range("A40").Value = 'knee pads airsoft'
$ws.Range("A41").Value = 'knee pads army'
$ws.Range("A42").Value = 'knee pads enduro'
$ws.Range("A43").Value = 'knee pads bike'
$ws.Range("A44").Value = 'knee pads caving'
$ws.Range("A45").Value = 'knee pads for women'
$ws.Range("A46").Value = 'knee pads gloves'
$ws.Range("A47").Value = 'knee pads green'
$ws.Range("A48").Value = 'knee pads kali'
$ws.Range("A49").Value = 'knee pads kuangmi'
$ws.Range("A50").Value = 'knee pads longboard'
$ws.Range("A51").Value = 'knee pads neoprene'
$ws.Range("A52").Value = 'knee pads ocp'
$ws.Range("A53").Value = 'knee pads orange'
$ws.Range("A54").Value = 'knee pads over pants'
$ws.Range("A55").Value = 'knee pads plastic'
$ws.Range("A56").Value = 'knee pads purple'
$ws.Range("A57").Value = 'knee pads razor'
$ws.Range("A58").Value = 'knee pads red'
$ws.Range("A59").Value = 'knee pads rollerblading'
$ws.Range("A60").Value = 'knee pads sailing'
$ws.Range("A61").Value = 'knee pads scooter'
$ws.Range("A62").Value = 'knee pads set'
$ws.Range("A63").Value = 'knee pads shooting'
$ws.Range("A64").Value = 'knee pads swat'
$ws.Range("A65").Value = 'knee pads teen'
$ws.Range("A66").Value = 'knee pads tsg'
$ws.Range("A67").Value = 'knee pads usmc'
$ws.Range("A68").Value = 'knee pads viper'
$ws.Range("A69").Value = 'knee pads white'
$ws.Range("A70").Value = 'knee pads yellow'
$ws.Range("A71").Value = 'mens leggings xs'
$ws.Range("A72").Value = 'xtextile compression pants men'
$ws.Range("A73").Value = 'crx men''s tights'
$ws.Range("A74").Value = 'yoga capri pants'
$ws.Range("A75").Value = 'dodoing kneepads'
$ws.Range("A76").Value = 'cavaliers basketball leggings'
$ws.Range("A77").Value = 'basketball knee pads kids'
$ws.Range("A78").Value = 'knee pads for toddlers'
$ws.Range("A79").Value = 'knee pads skating'
$ws.Range("A80").Value = 'knee pads skateboarding'
$ws.Range("A81").Value = 'basketball knee pads kids boys'
$ws.Range("A82").Value = 'knee pads for dance'
$ws.Range("A83").Value = 'knee pads rollerblade'
$ws.Range("A84").Value = 'knee pads tan'
$ws.Range("A85").Value = 'knee pad dancer'
$ws.Range("A86").Value = 'knee pads adidas'
$ws.Range("A87").Value = 'knee pads basketball mcdavid'
$ws.Range("A88").Value = 'knee pads dancing'
$ws.Range("A89").Value = 'knee pads dodgeball'
$ws.Range("A90").Value = 'knee pads pair'
$ws.Range("A91").Value = 'knee pads longboarding'
$ws.Range("A92").Value = 'knee pads nba'
$ws.Range("A93").Value = 'knee pads pole'
$ws.Range("A94").Value = 'knee pad and elbow pads'
$ws.Range("A95").Value = 'knee pad adidas'
$ws.Range("A96").Value = 'knee pad asics'
$ws.Range("A97").Value = 'knee pad for kids'
$ws.Range("A98").Value = 'knee pad military'
$ws.Range("A99").Value = 'knee pad mma'
$ws.Range("A100").Value = 'knee pad mizuno'
